# Apply the updated crypto symbol list values (price, volume%, hour)
# scraped by the GitHub Actions job. Every cell keeps its original
# inline-string (text) type, matching how the source feed writes
# these columns, so the write must not let Excel auto-coerce a
# numeric-looking token (e.g. "307.82", "9") into a real number.
#
# Recipe per cell: force the number format to Text ("@") so the
# COM layer stores the literal string, write the value, then
# ClearFormats() to drop the transient Text number-format again so
# the cell's style index is left exactly as it was before (style 0,
# i.e. no explicit <c s="..."> attribute) -- only the cell content
# changes, nothing else.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = [ordered]@{
    'D2' = '307.82'
    'G2' = '9'
    'D3' = '40.96'
    'E3' = '0.48%'
    'G3' = '9'
    'D4' = '5.207'
    'E4' = '1.79%'
    'G4' = '9'
    'D5' = '0.07667'
    'E5' = '0.63%'
    'G5' = '9'
    'D6' = '1.630'
    'E6' = '1.37%'
    'G6' = '9'
    'D7' = '0.9154'
    'E7' = '1.25%'
    'G7' = '9'
    'E8' = '-2.48%'
    'G8' = '9'
    'E9' = '9.48%'
    'G9' = '9'
    'D10' = '0.1822'
    'E10' = '2.22%'
    'G10' = '9'
    'D11' = '0.09125'
    'E11' = '-1.13%'
    'G11' = '9'
    'D12' = '0.04205'
    'E12' = '1.19%'
    'G12' = '9'
    'D13' = '0.1051'
    'E13' = '-0.17%'
    'G13' = '9'
    'D14' = '0.001261'
    'E14' = '0.08%'
    'G14' = '9'
    'D15' = '0.005728'
    'E15' = '-2.43%'
    'G15' = '9'
    'G16' = '9'
    'E17' = '-0.40%'
    'G17' = '9'
    'D18' = '4.303'
    'E18' = '1.24%'
    'G18' = '9'
    'D19' = '0.3336'
    'G19' = '9'
    'D20' = '7.354'
    'E20' = '12.55%'
    'G20' = '9'
    'D21' = '0.1383'
    'E21' = '1.35%'
    'G21' = '9'
    'E22' = '-1.67%'
    'G22' = '9'
    'D23' = '0.04024'
    'E23' = '-1.01%'
    'G23' = '9'
    'D24' = '0.001262'
    'E24' = '2.60%'
    'G24' = '9'
    'D25' = '0.004273'
    'E25' = '3.69%'
    'G25' = '9'
    'D26' = '0.0001302'
    'E26' = '0.06%'
    'G26' = '9'
    'G27' = '9'
    'G28' = '9'
    'G29' = '9'
    'G30' = '9'
    'G31' = '9'
    'G32' = '9'
    'G33' = '9'
    'G34' = '9'
    'G35' = '9'
    'G36' = '9'
    'G37' = '9'
    'E38' = '3.17%'
    'G38' = '9'
    'D39' = '0.05346'
    'E39' = '2.99%'
    'G39' = '9'
    'D40' = '0.007840'
    'E40' = '1.25%'
    'G40' = '9'
    'D41' = '0.1312'
    'E41' = '0.61%'
    'G41' = '9'
    'D42' = '0.006511'
    'E42' = '-7.11%'
    'G42' = '9'
    'D43' = '0.001913'
    'E43' = '-1.92%'
    'G43' = '9'
    'D44' = '0.008247'
    'E44' = '-6.09%'
    'G44' = '9'
    'D45' = '0.3332'
    'E45' = '-0.02%'
    'G45' = '9'
    'D46' = '0.00006718'
    'E46' = '-3.09%'
    'G46' = '9'
    'E47' = '0.15%'
    'G47' = '9'
    'D48' = '0.3816'
    'E48' = '1,123.97%'
    'G48' = '9'
    'D49' = '0.003105'
    'E49' = '-26.08%'
    'G49' = '9'
    'D50' = '0.00002103'
    'E50' = '0.15%'
    'G50' = '9'
    'E51' = '0.15%'
    'G51' = '9'
}

foreach ($cellRef in $updates.Keys) {
    $newValue = $updates[$cellRef]
    $cell = $ws.Range($cellRef)
    $cell.NumberFormat = "@"
    $cell.Value = $newValue
    $cell.ClearFormats()
}

